# Add a new "2025" worksheet (Negociado de Policia regional data) after the
# existing "2024" sheet, populate it with the new year's figures, and make
# it the active/selected sheet with the same selection the author left
# behind (B2:E14) when the workbook was saved.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2025"

# Header row (column titles reused from the other yearly sheets)
$headers = @("Región", "Víctimas: Mujeres", "Víctimas: Hombres", "Ofensores: Mujeres", "Ofensores: Hombres")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Regional data rows for 2025
$data = @(
    @("San Juan", 212, 88, 72, 215),
    @("Arecibo", 99, 41, 24, 118),
    @("Ponce", 119, 44, 30, 127),
    @("Humacao", 46, 8, 4, 50),
    @("Mayaguez", 51, 13, 2, 58),
    @("Caguas", 139, 61, 48, 151),
    @("Bayamón", 243, 70, 34, 280),
    @("Carolina", 112, 46, 24, 132),
    @("Guayama", 56, 18, 24, 49),
    @("Aguadilla", 75, 29, 28, 75),
    @("Utuado", 40, 20, 17, 37),
    @("Fajardo", 29, 15, 4, 35),
    @("Aibonito", 42, 17, 9, 54)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newSheet.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Match column widths (best-fit) of the sibling yearly sheets as closely as
# the engine's column-width quantization allows
$newSheet.Columns.Item(1).ColumnWidth = 9
$newSheet.Columns.Item(2).ColumnWidth = 16
$newSheet.Columns.Item(3).ColumnWidth = 16.833333333333334
$newSheet.Columns.Item(4).ColumnWidth = 17.833333333333334
$newSheet.Columns.Item(5).ColumnWidth = 18.666666666666668

# Make the new sheet the active tab with the same range selected as in the
# saved workbook
$newSheet.Activate()
$newSheet.Range("B2:E14").Select() | Out-Null
